{"js": "const pairs = [\n  [\"2025-10-22 Wednesday\", \"2025-10-23 Thursday\"],\n  [\"39\u00d785=3315\", \"58\u00d743=2494\"],\n  [\"66\u00d792=6072\", \"98\u00d759=5782\"],\n  [\"45\u00d793=4185\", \"85\u00d727=2295\"],\n  [\"81\u00d740=3240\", \"37\u00d720=740\"],\n  [\"64\u00d713=832\", \"69\u00d711=759\"],\n  [\"64\u00d739=2496\", \"25\u00d739=975\"],\n  [\"30\u00d742=1260\", \"84\u00d726=2184\"],\n  [\"93\u00d726=2418\", \"58\u00d759=3422\"],\n  [\"75\u00d764=4800\", \"24\u00d799=2376\"],\n  [\"42\u00d761=2562\", \"86\u00d717=1462\"],\n  [\"81\u00d731=2511\", \"73\u00d790=6570\"],\n  [\"22\u00d725=550\", \"51\u00d715=765\"],\n  [\"89\u00d746=4094\", \"79\u00d791=7189\"],\n  [\"45\u00d730=1350\", \"29\u00d726=754\"],\n  [\"31\u00d774=2294\", \"34\u00d750=1700\"],\n  [\"12\u00d758=696\", \"24\u00d712=288\"],\n  [\"96\u00d762=5952\", \"62\u00d791=5642\"],\n  [\"60\u00d746=2760\", \"68\u00d759=4012\"],\n  [\"74\u00d728=2072\", \"78\u00d779=6162\"],\n  [\"35\u00d725=875\", \"88\u00d741=3608\"],\n  [\"69\u00d714=966\", \"98\u00d742=4116\"],\n  [\"22\u00d769=1518\", \"84\u00d739=3276\"],\n  [\"63\u00d765=4095\", \"41\u00d791=3731\"],\n  [\"35\u00d787=3045\", \"63\u00d764=4032\"],\n  [\"56\u00d776=4256\", \"96\u00d748=4608\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"2025-10-22 Wednesday\", \"2025-10-23 Thursday\")\n  ,@(\"39\u00d785=3315\", \"58\u00d743=2494\")\n  ,@(\"66\u00d792=6072\", \"98\u00d759=5782\")\n  ,@(\"45\u00d793=4185\", \"85\u00d727=2295\")\n  ,@(\"81\u00d740=3240\", \"37\u00d720=740\")\n  ,@(\"64\u00d713=832\", \"69\u00d711=759\")\n  ,@(\"64\u00d739=2496\", \"25\u00d739=975\")\n  ,@(\"30\u00d742=1260\", \"84\u00d726=2184\")\n  ,@(\"93\u00d726=2418\", \"58\u00d759=3422\")\n  ,@(\"75\u00d764=4800\", \"24\u00d799=2376\")\n  ,@(\"42\u00d761=2562\", \"86\u00d717=1462\")\n  ,@(\"81\u00d731=2511\", \"73\u00d790=6570\")\n  ,@(\"22\u00d725=550\", \"51\u00d715=765\")\n  ,@(\"89\u00d746=4094\", \"79\u00d791=7189\")\n  ,@(\"45\u00d730=1350\", \"29\u00d726=754\")\n  ,@(\"31\u00d774=2294\", \"34\u00d750=1700\")\n  ,@(\"12\u00d758=696\", \"24\u00d712=288\")\n  ,@(\"96\u00d762=5952\", \"62\u00d791=5642\")\n  ,@(\"60\u00d746=2760\", \"68\u00d759=4012\")\n  ,@(\"74\u00d728=2072\", \"78\u00d779=6162\")\n  ,@(\"35\u00d725=875\", \"88\u00d741=3608\")\n  ,@(\"69\u00d714=966\", \"98\u00d742=4116\")\n  ,@(\"22\u00d769=1518\", \"84\u00d739=3276\")\n  ,@(\"63\u00d765=4095\", \"41\u00d791=3731\")\n  ,@(\"35\u00d787=3045\", \"63\u00d764=4032\")\n  ,@(\"56\u00d776=4256\", \"96\u00d748=4608\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}"}
